$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 228, pushing the existing rows
# (old 228..284) down to 230..286.
$ws.Rows("228:229").Insert()

# Row 228: new "Primera" quote for Provincia de Quillota
$ws.Cells.Item(228, 1).Value = 4
$ws.Cells.Item(228, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(228, 3).Value = "Los Lagos"
$ws.Cells.Item(228, 4).Value = 44551
$ws.Cells.Item(228, 5).Value = 10
$ws.Cells.Item(228, 6).Value = "Fruta"
$ws.Cells.Item(228, 7).Value = 100106
$ws.Cells.Item(228, 8).Value = "Oleaginosos"
$ws.Cells.Item(228, 9).Value = 100106002
$ws.Cells.Item(228, 10).Value = "Palta"
$ws.Cells.Item(228, 11).Value = "Hass"
$ws.Cells.Item(228, 12).Value = "Primera"
$ws.Cells.Item(228, 13).Value = 400
$ws.Cells.Item(228, 14).Value = 3900
$ws.Cells.Item(228, 15).Value = 4000
$ws.Cells.Item(228, 16).Value = 3950
$ws.Cells.Item(228, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(228, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(228, 19).Value = 3950
$ws.Cells.Item(228, 20).Value = 1

# Row 229: new "Segunda" quote for Provincia de Quillota
$ws.Cells.Item(229, 1).Value = 4
$ws.Cells.Item(229, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(229, 3).Value = "Los Lagos"
$ws.Cells.Item(229, 4).Value = 44551
$ws.Cells.Item(229, 5).Value = 10
$ws.Cells.Item(229, 6).Value = "Fruta"
$ws.Cells.Item(229, 7).Value = 100106
$ws.Cells.Item(229, 8).Value = "Oleaginosos"
$ws.Cells.Item(229, 9).Value = 100106002
$ws.Cells.Item(229, 10).Value = "Palta"
$ws.Cells.Item(229, 11).Value = "Hass"
$ws.Cells.Item(229, 12).Value = "Segunda"
$ws.Cells.Item(229, 13).Value = 200
$ws.Cells.Item(229, 14).Value = 3500
$ws.Cells.Item(229, 15).Value = 3500
$ws.Cells.Item(229, 16).Value = 3500
$ws.Cells.Item(229, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(229, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(229, 19).Value = 3500
$ws.Cells.Item(229, 20).Value = 1

# Make sure the date cells keep the same numeric date format used
# throughout column D (style index 2 -> yyyy-mm-dd hh:mm:ss).
$ws.Range("D228:D229").NumberFormat = $ws.Range("D230").NumberFormat
